# Update the fitted NRP1 excel sheet data
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: update raw kon (G7:G11) values for the "Lu et al. 2023 / VEGF165:VEGFR2" block ---
$sheet2.Range("G7").Value = 128404400
$sheet2.Range("G8").Value = 91080760
$sheet2.Range("G9").Value = 32735990
$sheet2.Range("G10").Value = 18411490
$sheet2.Range("G11").Value = 6732888

# Recalculate so dependent formulas (H7, K7:K11, L7 on Sheet2 and G3, K3 on Sheet1) update
$excel.Calculate()

# --- Update selections to match the saved view state captured in the diff ---
$sheet1.Activate()
$sheet1.Range("E25").Select()

$sheet2.Activate()
$sheet2.Range("G23").Select()

# Restore the originally active sheet (Sheet1 has tabSelected in the workbook)
$sheet1.Activate()
